# Auto-generated edit script applying F-column (想去人数/want-to-go count) corrections
# and one G-column (最低票价/min ticket price) availability-status change,
# per the commit diff (gh-pages data refresh at 456a3b4).
#
# Sheet index -> name mapping (1-based, matches $wb.Worksheets.Item(n)):
#   1 = 展览 (Exhibitions)
#   2 = 演出 (Performances)
#   3 = 本地生活 (Local Life)
#   4 = 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

# --- Sheet 1 (展览) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 6816
$ws.Range("F3").Value = 831
$ws.Range("F4").Value = 147
$ws.Range("F6").Value = 749
$ws.Range("F7").Value = 749
$ws.Range("F11").Value = 1117
$ws.Range("F12").Value = 877
$ws.Range("F13").Value = 15
$ws.Range("F14").Value = 708
$ws.Range("F15").Value = 1020
$ws.Range("F16").Value = 1364
$ws.Range("F17").Value = 49
$ws.Range("F18").Value = 125
$ws.Range("F19").Value = 548
$ws.Range("F21").Value = 582
$ws.Range("F25").Value = 1070
$ws.Range("F27").Value = 731
$ws.Range("F28").Value = 559
$ws.Range("F29").Value = 473
$ws.Range("F30").Value = 460
$ws.Range("F31").Value = 96
$ws.Range("F32").Value = 1011
$ws.Range("F33").Value = 1136
$ws.Range("F34").Value = 278
$ws.Range("F35").Value = 2375
$ws.Range("F37").Value = 1287
$ws.Range("F38").Value = 449
$ws.Range("F40").Value = 3899
# --- Sheet 2 (演出) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("G2").Value = "不可售"
$ws.Range("F8").Value = 168
$ws.Range("F15").Value = 339
$ws.Range("F22").Value = 243
$ws.Range("F23").Value = 248
$ws.Range("F30").Value = 1706
# --- Sheet 3 (本地生活) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 1260
$ws.Range("F5").Value = 1647
$ws.Range("F8").Value = 976
# --- Sheet 4 (全部类型) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1260
$ws.Range("F4").Value = 1647
$ws.Range("F7").Value = 976
$ws.Range("F8").Value = 6816
$ws.Range("F9").Value = 831
$ws.Range("F11").Value = 147
$ws.Range("F13").Value = 749
$ws.Range("F14").Value = 749
$ws.Range("F17").Value = 1117
$ws.Range("F18").Value = 877
$ws.Range("F19").Value = 708
$ws.Range("F20").Value = 168
$ws.Range("F21").Value = 168
$ws.Range("F23").Value = 1020
$ws.Range("F24").Value = 1365
$ws.Range("F25").Value = 49
$ws.Range("F26").Value = 125
$ws.Range("F27").Value = 548
$ws.Range("F29").Value = 582
$ws.Range("F31").Value = 339
$ws.Range("F33").Value = 1070
$ws.Range("F35").Value = 731
$ws.Range("F36").Value = 559
$ws.Range("F37").Value = 473
$ws.Range("F38").Value = 460
$ws.Range("F39").Value = 96
$ws.Range("F41").Value = 248
$ws.Range("F42").Value = 1011
$ws.Range("F43").Value = 1136
$ws.Range("F44").Value = 278
$ws.Range("F45").Value = 2375
$ws.Range("F47").Value = 1706
$ws.Range("F48").Value = 1706
$ws.Range("F49").Value = 1287
$ws.Range("F50").Value = 449
$ws.Range("F51").Value = 3899
